$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (rows 2-351). Bump it from 45171 (2023-09-02) to 45172 (2023-09-03)
# for every row, keeping the existing date formatting/style intact.
$ws.Range("C2:C351").Value = 45172
